$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price / 1h-volume-change data (scheduled GitHub Actions run).
# Rows 16/17 additionally swap Coin/Link/Price/Volume (ShibaInu <-> WrappedEther rank order).
# Force text format on Price column so values like "1.00" / "539.48" are stored as
# literal strings (matching the sheet's existing inline-string convention) rather than numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.402.54'
$ws.Range("E2").Value = '  +0.09%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.605.47'
$ws.Range("E3").Value = '  +0.42%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.48'
$ws.Range("E5").Value = '  +3.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.07'
$ws.Range("E6").Value = '  +0.38%  '

$ws.Range("E7").Value = '  +0.30%  '

$ws.Range("E8").Value = '  +0.31%  '

$ws.Range("E9").Value = '  -1.11%  '

$ws.Range("E10").Value = '  +1.61%  '

$ws.Range("E11").Value = '  +0.91%  '

$ws.Range("E12").Value = '  +1.85%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.065.98'
$ws.Range("E13").Value = '  +0.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '59.298.31'
$ws.Range("E14").Value = '  +0.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.57'
$ws.Range("E15").Value = '  +0.81%  '

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.610.23'
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000134'
$ws.Range("E17").Value = '  +0.74%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '342.96'
$ws.Range("E18").Value = '  +1.34%  '

$ws.Range("E19").Value = '  +0.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.12'
$ws.Range("E20").Value = '  -0.70%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.42'
$ws.Range("E21").Value = '  -1.17%  '

$ws.Range("E22").Value = '  +0.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.44'
$ws.Range("E23").Value = '  +1.55%  '

$ws.Range("E24").Value = '  -0.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.409'
$ws.Range("E25").Value = '  +1.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.46%  '

$ws.Range("E27").Value = '  +1.84%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0739'
$ws.Range("E29").Value = '  +1.88%  '

$ws.Range("E30").Value = '  +5.99%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.84'
$ws.Range("E31").Value = '  -1.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.78'
$ws.Range("E32").Value = '  -0.18%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.32'
$ws.Range("E33").Value = '  +0.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.98'
$ws.Range("E34").Value = '  -0.54%  '

$ws.Range("E35").Value = '  -1.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '36.98'
$ws.Range("E36").Value = '  +1.68%  '

$ws.Range("E37").Value = '  +0.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.835'
$ws.Range("E38").Value = '  +0.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.825'
$ws.Range("E39").Value = '  +0.18%  '

$ws.Range("E40").Value = '  +0.39%  '

$ws.Range("E41").Value = '  +0.41%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '274.16'
$ws.Range("E42").Value = '  -0.55%  '

$ws.Range("E43").Value = '  +0.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.596'
$ws.Range("E44").Value = '  +0.83%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0957'
$ws.Range("E45").Value = '  +0.42%  '

$ws.Range("E46").Value = '  +0.62%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.947.04'
$ws.Range("E47").Value = '  -1.70%  '

$ws.Range("E48").Value = '  +0.95%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.28'
$ws.Range("E49").Value = '  +0.91%  '

$ws.Range("E50").Value = '  -2.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.96'
$ws.Range("E51").Value = '  -2.15%  '
